$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style for unstyled data cells (used to strip any quote-prefix
# style Excel applies when a numeric-looking string is force-typed as text).
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "71.333.90"
$ws.Range("E2").Value = "  +6.91%  "
$ws.Range("D3").Value = "3.688.79"
$ws.Range("E3").Value = "  +19.40%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'599.07"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  +3.89%  "
$ws.Range("D6").Value = "'184.18"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  +6.77%  "
$ws.Range("D7").Value = "3.687.44"
$ws.Range("E7").Value = "  +19.47%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +4.25%  "
$ws.Range("E10").Value = "  +8.32%  "
$ws.Range("D11").Value = "'6.62"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = "  +4.09%  "
$ws.Range("D12").Value = "'0.499"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  +5.77%  "
$ws.Range("D13").Value = "'39.85"
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = "  +11.60%  "
$ws.Range("E14").Value = "  +6.29%  "
$ws.Range("D15").Value = "4.302.72"
$ws.Range("E15").Value = "  +19.37%  "
$ws.Range("D16").Value = "3.682.53"
$ws.Range("E16").Value = "  +19.17%  "
$ws.Range("D17").Value = "71.264.09"
$ws.Range("E17").Value = "  +6.95%  "
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("E19").Value = "  +7.57%  "
$ws.Range("D20").Value = "'16.88"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "'515.44"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  +6.61%  "
$ws.Range("D22").Value = "'9.19"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = "  +18.71%  "
$ws.Range("E23").Value = "  +8.32%  "
$ws.Range("D24").Value = "'87.46"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  +4.99%  "
$ws.Range("E25").Value = "  +8.52%  "
$ws.Range("E26").Value = "  +5.75%  "
$ws.Range("D27").Value = "'10.89"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  +8.31%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  +12.66%  "
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("D31").Value = "'31.71"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "  +13.57%  "
$ws.Range("E32").Value = "  +7.19%  "
$ws.Range("E33").Value = "  +18.00%  "
$ws.Range("E34").Value = "  +3.93%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'6.12"
$ws.Range("D36").Style = $plainStyle
$ws.Range("E37").Value = "  +8.30%  "
$ws.Range("D38").Value = "'0.341"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  +10.29%  "
$ws.Range("E39").Value = "  +9.70%  "
$ws.Range("D40").Value = "'50.88"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = "  +3.87%  "
$ws.Range("E41").Value = "  +4.18%  "
$ws.Range("D42").Value = "'45.13"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  -6.37%  "
$ws.Range("D43").Value = "3.170.78"
$ws.Range("E43").Value = "  +14.18%  "
$ws.Range("D44").Value = "'8.78"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  +6.11%  "
$ws.Range("D45").Value = "'409.68"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = "  +11.50%  "
$ws.Range("E46").Value = "  +6.24%  "
$ws.Range("E47").Value = "  +6.25%  "
$ws.Range("D48").Value = "'28.13"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  +15.25%  "
$ws.Range("D49").Value = "'136.83"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = "  +1.96%  "
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("E51").Value = "  +12.72%  "
